# Updated symbol list on Sun Feb 12 20:51:57 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for most rows, and for rows
# 14-22 shifts the Coin/Link/Price rows down by one (BitForexToken moves to
# the top of that block) while updating Price/Volume for the shifted rows.
#
# All D/E values are plain text in the source workbook (no numeric/percent
# cell format applied), so we write them with a leading "'" (quote-prefix)
# to force Excel to keep them as text instead of auto-converting to
# number/percentage, then reset the cell style back to "Normal" so the
# quote-prefix flag doesn't linger as a style change on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'320.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'3.66%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'41.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.24%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.258"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.64%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07744"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.56%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.758"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'9.40%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.9447"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.93%"
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'-0.82%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.64%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.46%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.09212"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.48%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.04152"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.33%"
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'0.47%"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001284"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.82%"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005825"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.71%"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007491"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1,897.31%"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.351"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.02%"
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.335"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.00%"
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3358"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.21%"
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'8.416"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'21.42%"
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1354"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.82%"
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2825"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.18%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04035"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.32%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001268"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.18%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004119"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.76%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001272"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D38").Value = "'0.02553"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'5.59%"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.05355"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.32%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.007774"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.78%"
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'1.20%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.007023"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.30%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.001992"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.98%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.008305"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.83%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.3452"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.22%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006691"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.63%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'0.00%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.1991"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'57.84%"
$ws.Range("E48").Style = "Normal"

$ws.Range("E50").Value = "'0.00%"
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'0.00%"
$ws.Range("E51").Style = "Normal"
